$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-01-23 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-01-24 Wednesday", 2)

# Update the division-problem answers in the table, cell by cell
# (several original strings repeat, so Find/Replace-all would be unsafe;
# addressing cells directly by row/column guarantees the correct mapping)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "40÷9=4, 4"
$t.Cell(1, 2).Range.Text = "36÷7=5, 1"
$t.Cell(1, 3).Range.Text = "19÷5=3, 4"
$t.Cell(1, 4).Range.Text = "21÷8=2, 5"
$t.Cell(1, 5).Range.Text = "82÷2=41, 0"

$t.Cell(5, 1).Range.Text = "74÷4=18, 2"
$t.Cell(5, 2).Range.Text = "89÷9=9, 8"
$t.Cell(5, 3).Range.Text = "12÷5=2, 2"
$t.Cell(5, 4).Range.Text = "71÷3=23, 2"
$t.Cell(5, 5).Range.Text = "12÷6=2, 0"

$t.Cell(9, 1).Range.Text = "59÷9=6, 5"
$t.Cell(9, 2).Range.Text = "61÷8=7, 5"
$t.Cell(9, 3).Range.Text = "28÷8=3, 4"
$t.Cell(9, 4).Range.Text = "82÷6=13, 4"
$t.Cell(9, 5).Range.Text = "36÷6=6, 0"

$t.Cell(13, 1).Range.Text = "13÷2=6, 1"
$t.Cell(13, 2).Range.Text = "74÷9=8, 2"
$t.Cell(13, 3).Range.Text = "24÷4=6, 0"
$t.Cell(13, 4).Range.Text = "74÷8=9, 2"
$t.Cell(13, 5).Range.Text = "86÷2=43, 0"

$t.Cell(17, 1).Range.Text = "85÷8=10, 5"
$t.Cell(17, 2).Range.Text = "85÷4=21, 1"
$t.Cell(17, 3).Range.Text = "88÷2=44, 0"
$t.Cell(17, 4).Range.Text = "85÷7=12, 1"
$t.Cell(17, 5).Range.Text = "18÷7=2, 4"
